# This script applies updated TPM-derived values to the NATMI
# ligand-receptor pair output sheet (Apoe-Ldlr), per commit
# "update scripts wuth new tpm".
#
# The underlying analysis was rerun with new TPM values, which changes
# the Ligand/Receptor average & total expression values (columns G, H,
# M, N), their derived-specificity counterparts (I, J, O, P), and the
# Edge expression weights / derived specificities that are computed
# from them (Q, R, S, T). Columns A-F, K, L (cluster/gene identifiers,
# counts, detection rates) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 35.42516366666666
$ws.Cells.Item(2, 8).Value = 106.275491
$ws.Cells.Item(2, 9).Value = 0.00832770193000585
$ws.Cells.Item(2, 10).Value = 0.008327701930005852
$ws.Cells.Item(2, 13).Value = 1.711937666666667
$ws.Cells.Item(2, 14).Value = 5.135813
$ws.Cells.Item(2, 15).Value = 0.09827532014408574
$ws.Cells.Item(2, 16).Value = 0.09827532014408574
$ws.Cells.Item(2, 17).Value = 60.6456720287981
$ws.Cells.Item(2, 18).Value = 545.8110482591829
$ws.Cells.Item(2, 19).Value = 0.0008184075732358457
$ws.Cells.Item(2, 20).Value = 0.0008184075732358458
$ws.Cells.Item(3, 7).Value = 35.42516366666666
$ws.Cells.Item(3, 8).Value = 106.275491
$ws.Cells.Item(3, 9).Value = 0.00832770193000585
$ws.Cells.Item(3, 10).Value = 0.008327701930005852
$ws.Cells.Item(3, 15).Value = 0.3329367223581701
$ws.Cells.Item(3, 16).Value = 0.3329367223581701
$ws.Cells.Item(3, 17).Value = 205.4551564001362
$ws.Cells.Item(3, 18).Value = 1849.096407601226
$ws.Cells.Item(3, 19).Value = 0.002772597785351955
$ws.Cells.Item(3, 20).Value = 0.002772597785351955
$ws.Cells.Item(4, 7).Value = 35.42516366666666
$ws.Cells.Item(4, 8).Value = 106.275491
$ws.Cells.Item(4, 9).Value = 0.00832770193000585
$ws.Cells.Item(4, 10).Value = 0.008327701930005852
$ws.Cells.Item(4, 13).Value = 3.605537
$ws.Cells.Item(4, 14).Value = 10.816611
$ws.Cells.Item(4, 15).Value = 0.2069790915087912
$ws.Cells.Item(4, 16).Value = 0.2069790915087912
$ws.Cells.Item(4, 17).Value = 127.7267383312223
$ws.Cells.Item(4, 18).Value = 1149.540644981001
$ws.Cells.Item(4, 19).Value = 0.001723660179828618
$ws.Cells.Item(4, 20).Value = 0.001723660179828618
$ws.Cells.Item(5, 7).Value = 35.42516366666666
$ws.Cells.Item(5, 8).Value = 106.275491
$ws.Cells.Item(5, 9).Value = 0.00832770193000585
$ws.Cells.Item(5, 10).Value = 0.008327701930005852
$ws.Cells.Item(5, 13).Value = 6.302642666666666
$ws.Cells.Item(5, 14).Value = 18.907928
$ws.Cells.Item(5, 15).Value = 0.3618088659889531
$ws.Cells.Item(5, 16).Value = 0.361808865988953
$ws.Cells.Item(5, 17).Value = 223.2721479991831
$ws.Cells.Item(5, 18).Value = 2009.449331992648
$ws.Cells.Item(5, 19).Value = 0.003013036391589433
$ws.Cells.Item(5, 20).Value = 0.003013036391589433
$ws.Cells.Item(6, 9).Value = 0.01070182047907406
$ws.Cells.Item(6, 10).Value = 0.01070182047907406
$ws.Cells.Item(6, 13).Value = 1.711937666666667
$ws.Cells.Item(6, 14).Value = 5.135813
$ws.Cells.Item(6, 15).Value = 0.09827532014408574
$ws.Cells.Item(6, 16).Value = 0.09827532014408574
$ws.Cells.Item(6, 17).Value = 77.93495736758966
$ws.Cells.Item(6, 18).Value = 701.414616308307
$ws.Cells.Item(6, 19).Value = 0.001051724833705536
$ws.Cells.Item(6, 20).Value = 0.001051724833705536
$ws.Cells.Item(7, 9).Value = 0.01070182047907406
$ws.Cells.Item(7, 10).Value = 0.01070182047907406
$ws.Cells.Item(7, 15).Value = 0.3329367223581701
$ws.Cells.Item(7, 16).Value = 0.3329367223581701
$ws.Cells.Item(7, 19).Value = 0.003563029033568458
$ws.Cells.Item(7, 20).Value = 0.003563029033568458
$ws.Cells.Item(8, 9).Value = 0.01070182047907406
$ws.Cells.Item(8, 10).Value = 0.01070182047907406
$ws.Cells.Item(8, 13).Value = 3.605537
$ws.Cells.Item(8, 14).Value = 10.816611
$ws.Cells.Item(8, 15).Value = 0.2069790915087912
$ws.Cells.Item(8, 16).Value = 0.2069790915087912
$ws.Cells.Item(8, 17).Value = 164.139955474781
$ws.Cells.Item(8, 18).Value = 1477.259599273029
$ws.Cells.Item(8, 19).Value = 0.002215053080248925
$ws.Cells.Item(8, 20).Value = 0.002215053080248925
$ws.Cells.Item(9, 9).Value = 0.01070182047907406
$ws.Cells.Item(9, 10).Value = 0.01070182047907406
$ws.Cells.Item(9, 13).Value = 6.302642666666666
$ws.Cells.Item(9, 14).Value = 18.907928
$ws.Cells.Item(9, 15).Value = 0.3618088659889531
$ws.Cells.Item(9, 16).Value = 0.361808865988953
$ws.Cells.Item(9, 17).Value = 286.9241077487547
$ws.Cells.Item(9, 18).Value = 2582.316969738792
$ws.Cells.Item(9, 19).Value = 0.003872013531551138
$ws.Cells.Item(9, 20).Value = 0.003872013531551138
$ws.Cells.Item(10, 7).Value = 51.06824600000001
$ws.Cells.Item(10, 8).Value = 153.204738
$ws.Cells.Item(10, 9).Value = 0.01200505761322374
$ws.Cells.Item(10, 10).Value = 0.01200505761322374
$ws.Cells.Item(10, 13).Value = 1.711937666666667
$ws.Cells.Item(10, 14).Value = 5.135813
$ws.Cells.Item(10, 15).Value = 0.09827532014408574
$ws.Cells.Item(10, 16).Value = 0.09827532014408574
$ws.Cells.Item(10, 17).Value = 87.42565389799934
$ws.Cells.Item(10, 18).Value = 786.8308850819941
$ws.Cells.Item(10, 19).Value = 0.001179800880287757
$ws.Cells.Item(10, 20).Value = 0.001179800880287757
$ws.Cells.Item(11, 7).Value = 51.06824600000001
$ws.Cells.Item(11, 8).Value = 153.204738
$ws.Cells.Item(11, 9).Value = 0.01200505761322374
$ws.Cells.Item(11, 10).Value = 0.01200505761322374
$ws.Cells.Item(11, 15).Value = 0.3329367223581701
$ws.Cells.Item(11, 16).Value = 0.3329367223581701
$ws.Cells.Item(11, 17).Value = 296.1802680077187
$ws.Cells.Item(11, 18).Value = 2665.622412069468
$ws.Cells.Item(11, 19).Value = 0.003996924533467707
$ws.Cells.Item(11, 20).Value = 0.003996924533467708
$ws.Cells.Item(12, 7).Value = 51.06824600000001
$ws.Cells.Item(12, 8).Value = 153.204738
$ws.Cells.Item(12, 9).Value = 0.01200505761322374
$ws.Cells.Item(12, 10).Value = 0.01200505761322374
$ws.Cells.Item(12, 13).Value = 3.605537
$ws.Cells.Item(12, 14).Value = 10.816611
$ws.Cells.Item(12, 15).Value = 0.2069790915087912
$ws.Cells.Item(12, 16).Value = 0.2069790915087912
$ws.Cells.Item(12, 17).Value = 184.128450478102
$ws.Cells.Item(12, 18).Value = 1657.156054302918
$ws.Cells.Item(12, 19).Value = 0.002484795918295746
$ws.Cells.Item(12, 20).Value = 0.002484795918295747
$ws.Cells.Item(13, 7).Value = 51.06824600000001
$ws.Cells.Item(13, 8).Value = 153.204738
$ws.Cells.Item(13, 9).Value = 0.01200505761322374
$ws.Cells.Item(13, 10).Value = 0.01200505761322374
$ws.Cells.Item(13, 13).Value = 6.302642666666666
$ws.Cells.Item(13, 14).Value = 18.907928
$ws.Cells.Item(13, 15).Value = 0.3618088659889531
$ws.Cells.Item(13, 16).Value = 0.361808865988953
$ws.Cells.Item(13, 17).Value = 321.8649061514294
$ws.Cells.Item(13, 18).Value = 2896.784155362864
$ws.Cells.Item(13, 19).Value = 0.004343536281172528
$ws.Cells.Item(13, 20).Value = 0.004343536281172528
$ws.Cells.Item(14, 7).Value = 4121.876464666667
$ws.Cells.Item(14, 8).Value = 12365.629394
$ws.Cells.Item(14, 9).Value = 0.9689654199776964
$ws.Cells.Item(14, 10).Value = 0.9689654199776964
$ws.Cells.Item(14, 13).Value = 1.711937666666667
$ws.Cells.Item(14, 14).Value = 5.135813
$ws.Cells.Item(14, 15).Value = 0.09827532014408574
$ws.Cells.Item(14, 16).Value = 0.09827532014408574
$ws.Cells.Item(14, 17).Value = 7056.395577209702
$ws.Cells.Item(14, 18).Value = 63507.56019488732
$ws.Cells.Item(14, 19).Value = 0.09522538685685661
$ws.Cells.Item(14, 20).Value = 0.09522538685685661
$ws.Cells.Item(15, 7).Value = 4121.876464666667
$ws.Cells.Item(15, 8).Value = 12365.629394
$ws.Cells.Item(15, 9).Value = 0.9689654199776964
$ws.Cells.Item(15, 10).Value = 0.9689654199776964
$ws.Cells.Item(15, 15).Value = 0.3329367223581701
$ws.Cells.Item(15, 16).Value = 0.3329367223581701
$ws.Cells.Item(15, 17).Value = 23905.62769670377
$ws.Cells.Item(15, 18).Value = 215150.6492703339
$ws.Cells.Item(15, 19).Value = 0.3226041710057819
$ws.Cells.Item(15, 20).Value = 0.3226041710057819
$ws.Cells.Item(16, 7).Value = 4121.876464666667
$ws.Cells.Item(16, 8).Value = 12365.629394
$ws.Cells.Item(16, 9).Value = 0.9689654199776964
$ws.Cells.Item(16, 10).Value = 0.9689654199776964
$ws.Cells.Item(16, 13).Value = 3.605537
$ws.Cells.Item(16, 14).Value = 10.816611
$ws.Cells.Item(16, 15).Value = 0.2069790915087912
$ws.Cells.Item(16, 16).Value = 0.2069790915087912
$ws.Cells.Item(16, 17).Value = 14861.57810278486
$ws.Cells.Item(16, 18).Value = 133754.2029250637
$ws.Cells.Item(16, 19).Value = 0.2005555823304179
$ws.Cells.Item(16, 20).Value = 0.2005555823304179
$ws.Cells.Item(17, 7).Value = 4121.876464666667
$ws.Cells.Item(17, 8).Value = 12365.629394
$ws.Cells.Item(17, 9).Value = 0.9689654199776964
$ws.Cells.Item(17, 10).Value = 0.9689654199776964
$ws.Cells.Item(17, 13).Value = 6.302642666666666
$ws.Cells.Item(17, 14).Value = 18.907928
$ws.Cells.Item(17, 15).Value = 0.3618088659889531
$ws.Cells.Item(17, 16).Value = 0.361808865988953
$ws.Cells.Item(17, 17).Value = 25978.71447293729
$ws.Cells.Item(17, 18).Value = 233808.4302564356
$ws.Cells.Item(17, 19).Value = 0.35058027978464
$ws.Cells.Item(17, 20).Value = 0.3505802797846399
